$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.440.37"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.09%  "
$ws.Range("D3").Value = "'1.838.41"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.57%  "
$ws.Range("D4").Value = "'1.027"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +2.49%  "
$ws.Range("D5").Value = "'318.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.01%  "
$ws.Range("D6").Value = "'1.024"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.20%  "
$ws.Range("D7").Value = "'0.4363"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.08%  "
$ws.Range("E8").Value = "  +3.35%  "
$ws.Range("E9").Value = "  +3.30%  "
$ws.Range("D10").Value = "'0.8728"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.44%  "
$ws.Range("D11").Value = "'21.37"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.80%  "
$ws.Range("D12").Value = "'1.883.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.62%  "
$ws.Range("D13").Value = "'5.466"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.50%  "
$ws.Range("E14").Value = "  +3.79%  "
$ws.Range("D15").Value = "'0.07145"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.17%  "
$ws.Range("D16").Value = "'82.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.94%  "
$ws.Range("D17").Value = "'1.031"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.87%  "
$ws.Range("D18").Value = "'0.000008986"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.58%  "
$ws.Range("D19").Value = "'1.025"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.24%  "
$ws.Range("D20").Value = "'15.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.07%  "
$ws.Range("D21").Value = "'27.447.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.99%  "
$ws.Range("D22").Value = "'5.234"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.24%  "
$ws.Range("D23").Value = "'11.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.03%  "
$ws.Range("D24").Value = "'2.084.54"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.35%  "
$ws.Range("D25").Value = "'156.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.88%  "
$ws.Range("D26").Value = "'1.899"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.77%  "
$ws.Range("E27").Value = "  +3.37%  "
$ws.Range("D28").Value = "'5.236"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.47%  "
$ws.Range("D29").Value = "'1.926"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.66%  "
$ws.Range("D30").Value = "'116.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.74%  "
$ws.Range("D31").Value = "'0.09055"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.39%  "
$ws.Range("D32").Value = "'1.202"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.76%  "
$ws.Range("D33").Value = "'0.7599"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.70%  "
$ws.Range("D34").Value = "'4.481"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.75%  "
$ws.Range("D35").Value = "'2.868"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.84%  "
$ws.Range("D36").Value = "'1.027"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.51%  "
$ws.Range("D37").Value = "'1.146"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.90%  "
$ws.Range("E38").Value = "  +4.38%  "
$ws.Range("D39").Value = "'0.05250"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.34%  "
$ws.Range("D40").Value = "'0.5169"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.59%  "
$ws.Range("D41").Value = "'2.780"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.19%  "
$ws.Range("D42").Value = "'0.1662"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.42%  "
$ws.Range("D43").Value = "'6.558"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.84%  "
$ws.Range("D44").Value = "'8.487"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.75%  "
$ws.Range("D45").Value = "'108.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.99%  "
$ws.Range("D46").Value = "'10.52"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.96%  "
$ws.Range("D47").Value = "'1.028"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.68%  "
$ws.Range("D48").Value = "'1.683"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.25%  "
$ws.Range("D49").Value = "'0.4628"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.48%  "
$ws.Range("D50").Value = "'1.900"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +10.76%  "
$ws.Range("D51").Value = "'0.06310"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.24%  "
